$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: clear the stray answer-option text that was left in G7
$ws.Range("G7").ClearContents()

# Move the active selection to G14 (single cell, matching the saved view state)
$ws.Range("G14").Select()
